# Edit Saldo.xlsx per commit diff:
#  1. Change the "Gustavo" (account 004444605) balance from 50700 to 62600.
#  2. Remove the "Elaine" (account 008384472) row entirely.
#  3. Remove the "Joao" (account 008026930) row entirely.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Export")

# 1. Update Gustavo's balance (row 4, column C) from 50700 to 62600.
$ws.Cells.Item(4, 3).Value = 62600

# 2. Delete the Elaine row (account 008384472) - originally row 10.
$ws.Rows.Item(10).Delete()

# 3. Delete the Joao row (account 008026930).
# It was originally row 272; after removing the Elaine row above, it shifts to row 271.
$ws.Rows.Item(271).Delete()
